# ---------------------------------------------------------------------------
# Recreates the "saving wip before pulling" commit:
#  - Sheet2 gets populated with a static (values-only) copy of the "case 2"
#    sheet's frequency / PSD2 columns (E17:G83 -> A3:B69), with a small
#    "Static Graph" header block, and becomes the active sheet.
#  - The "case 2" sheet's scroll/selection state moves on to a new range.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("case 2")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Populate Sheet2 with the "Static Graph" header + values -------------
$ws2.Range("A1").Value = "Static Graph"
$ws2.Range("A2").Value = "freq (1/micron)"
$ws2.Range("B2").Value = "PSD2 (nm^4)"

$rows = 67
$data = New-Object 'object[,]' $rows,2
$src = @(
  @("1E-4","11548912363.159761"),
  @("1.2589254117941672E-4","11455926786.628166"),
  @("1.584893192461112E-4","11311158638.043528"),
  @("1.995262314968879E-4","11088045619.768627"),
  @("2.5118864315095817E-4","10749535127.983898"),
  @("3.1622776601683783E-4","10248050335.388151"),
  @("3.9810717055349746E-4","9531081019.0810699"),
  @("5.0118723362727253E-4","8557346111.9564781"),
  @("6.3095734448019407E-4","7325450882.548583"),
  @("7.943282347242824E-4","5904340093.7899733"),
  @("1.0000000000000011E-3","4437608810.8350315"),
  @("1.2589254117941697E-3","3098641729.0728121"),
  @("1.5848931924611165E-3","2017596708.5609155"),
  @("1.9952623149688833E-3","1236912840.2259886"),
  @("2.5118864315095868E-3","723006432.6861155"),
  @("3.1622776601683876E-3","407982232.27111661"),
  @("3.981071705534983E-3","224589621.51376194"),
  @("5.0118723362727394E-3","121578002.59917131"),
  @("6.3095734448019537E-3","65088924.178863503"),
  @("7.9432823472428398E-3","34596345.808077008"),
  @("1.000000000000004E-2","18303792.584701691"),
  @("1.2589254117941722E-2","9655322.0259411838"),
  @("1.5848931924611197E-2","5083645.6867233664"),
  @("1.9952623149688879E-2","2673410.6034739409"),
  @("2.5118864315095916E-2","1404844.1178447171"),
  @("3.1622776601683951E-2","737875.85696740448"),
  @("3.9810717055349922E-2","387442.71579651738"),
  @("5.0118723362727491E-2","203399.08889387906"),
  @("6.3095734448019664E-2","106767.31040998822"),
  @("7.9432823472428596E-2","56039.551812809295"),
  @("0.10000000000000055","29412.382884953644"),
  @("0.12589254117941742","15436.63542488416"),
  @("0.15848931924611223","8101.5255160702845"),
  @("0.17782794100389329","5869.1036907380994"),
  @("0.19952623149688914","4251.8284806322254"),
  @("0.22387211385683531","3080.1997251348484"),
  @("0.25118864315095951","2231.4202647181401"),
  @("0.28183829312644709","1616.5284587868664"),
  @("0.31622776601683988","1171.0755910713826"),
  @("0.35481338923357758","848.37171291825985"),
  @("0.3981071705534997","614.59238202784229"),
  @("0.44668359215096581","445.23363121604905"),
  @("0.50118723362727535","322.54371855481952"),
  @("0.56234132519035251","233.66253626639855"),
  @("0.63095734448019714","169.27369792497802"),
  @("0.70794578438414224","122.6280425651394"),
  @("0.72443596007499445","114.97140814719737"),
  @("0.74131024130092205","107.79283710357404"),
  @("0.75857757502918843","101.06248027486842"),
  @("0.7762471166286965","94.752352206728673"),
  @("0.79432823472428638","88.836214785582456"),
  @("0.81283051616410418","83.289468139735447"),
  @("0.83176377110267607","78.089048352226143"),
  @("0.85113803820238165","73.213331560138201"),
  @("0.87096358995608603","68.642044041629035"),
  @("0.89125093813375089","64.356177916829949"),
  @("0.9120108393559152","60.337912112103069"),
  @("0.93325430079699667","56.570538259044433"),
  @("0.95499258602144177","53.038391220116445"),
  @("0.97723722095581667","49.726783952047548"),
  @("0.98174794301999047","49.089700026990975"),
  @("0.98627948563121648","48.460778214952299"),
  @("0.9908319448927736","47.83991394571207"),
  @("0.99540541735153309","47.22700398876642"),
  @("0.99655208013477448","47.075007389670816"),
  @("0.99770006382255938","46.923499979971538"),
  @("0.9988493699365113","46.772480185251993")
)
for ($i = 0; $i -lt $rows; $i++) {
  $data[$i,0] = [double]$src[$i][0]
  $data[$i,1] = [double]$src[$i][1]
}
$ws2.Range("A3:B69").Value = $data

# --- View state: Sheet2 becomes the active/selected sheet ----------------
$ws2.Activate()
$ws2.Range("A3:B69").Select()

# --- "case 2" sheet: selection moves to G17:G83 ---------------------------
$ws1.Range("G17:G83").Select()

# Re-activate Sheet2 so it ends up as the tab shown when the file is opened.
$ws2.Activate()
